$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move "tunnit:" label + array total formula from H1/I1 to E1/F1 ---
$ws.Range("H1").Copy()
$ws.Range("E1").PasteSpecial(-4104)   # xlPasteAll
$ws.Range("H1").ClearContents()

$ws.Range("F1").FormulaArray = "=SUM(B2:B999/60)"
$ws.Range("I1").ClearContents()

# --- Add new diary rows 14 and 15 ---
$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14").Value = 44492
$ws.Range("B14").Value = 60
$ws.Range("C14").Value = "Debuggia päivämäärien ja scorejen kanssa"

$ws.Range("A2").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A15").Value = 44590
$ws.Range("B15").Value = 80
$ws.Range("C15").Value = "Päivämäärän hallinta toimii nyt, alotettu draw.io dokumentaatio"

# --- Wrap text + widen column C (applies to existing + new rows) ---
$ws.Columns("C").WrapText = $true
$ws.Columns("C").ColumnWidth = 76

# --- Update selection to match the new "next empty row" cell ---
$ws.Range("H15").Select()
